$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.797.66"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.054.04"
$ws.Range("E3").Value = "  -4.07%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.93"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.18"
$ws.Range("E6").Value = "  -4.16%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.052.85"
$ws.Range("E8").Value = "  -4.13%  "
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.26"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.441"
$ws.Range("E12").Value = "  -2.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000235"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.80"
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.560.02"
$ws.Range("E16").Value = "  -3.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.868.91"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.054.74"
$ws.Range("E18").Value = "  -3.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.39"
$ws.Range("E19").Value = "  -2.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "449.43"
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.56"
$ws.Range("E21").Value = "  -2.57%  "
$ws.Range("E22").Value = "  -4.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.36"
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.20"
$ws.Range("E24").Value = "  -2.90%  "
$ws.Range("E25").Value = "  -4.03%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  -4.19%  "
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.42"
$ws.Range("E30").Value = "  -4.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.47"
$ws.Range("E31").Value = "  -6.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.94"
$ws.Range("E32").Value = "  -5.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0977"
$ws.Range("E33").Value = "  -6.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.34"
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.974"
$ws.Range("E35").Value = "  -5.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.73"
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "50.44"
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0696"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  -2.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.95"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "382.58"
$ws.Range("E41").Value = "  -6.27%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.109"
$ws.Range("E42").Value = "  -2.46%  "
$ws.Range("E43").Value = "  -6.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.700.16"
$ws.Range("E44").Value = "  -5.94%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.16"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.02"
$ws.Range("E48").Value = "  -5.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.08"
$ws.Range("E49").Value = "  -6.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.108"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.08"
$ws.Range("E51").Value = "  -5.55%  "

Write-Host "Applied cryptos update"